$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Cursor-Werte" row (row 2). Excel shifts the rows below it
# (Print graph, Automatenmodus) up by one, matching the diff exactly.
$ws.Rows.Item(2).Delete()

# Reflect the resulting selection state recorded in the saved file.
$ws.Range("A2:A3").Select()
